$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# format them as Text first so they are stored as the literal string, then
# restore the Normal style so no stray number-format style remains on the cell.
$textCells = @{
    'D4' = '1.001'
    'D6' = '303.89'
    'D7' = '0.3786'
    'D8' = '52.14'
    'D9' = '0.3603'
    'D10' = '0.08068'
    'D11' = '1.221'
    'D12' = '1.001'
    'D13' = '22.55'
    'D14' = '6.532'
    'D15' = '0.00001243'
    'D16' = '7.198'
    'D18' = '93.49'
    'D19' = '0.06909'
    'D20' = '17.86'
    'D22' = '6.405'
    'D24' = '12.68'
    'D25' = '3.168'
    'D26' = '2.445'
    'D27' = '21.05'
    'D28' = '149.89'
    'D29' = '5.285'
    'D30' = '134.52'
    'D31' = '2.293'
    'D33' = '6.747'
    'D34' = '10.91'
    'D35' = '0.9451'
    'D36' = '0.02772'
    'D37' = '0.2511'
    'D38' = '0.08814'
    'D39' = '6.027'
    'D40' = '0.07096'
    'D41' = '1.355'
    'D42' = '0.7001'
    'D43' = '15.97'
    'D44' = '12.21'
    'D46' = '0.6403'
    'D47' = '2.300'
    'D48' = '3.986'
    'D49' = '0.07974'
    'D51' = '125.07'
}
foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}

# Remaining cells: plain text / links / percentages / multi-dot price strings
# that Excel will not misinterpret as numbers, so they can be set directly.
$plainCells = @{
    'D2' = '23.308.66'
    'E2' = '  -1.05%  '
    'D3' = '1.621.77'
    'E3' = '  -0.96%  '
    'E4' = '  +0.55%  '
    'E5' = '  +0.54%  '
    'E6' = '  -1.27%  '
    'E7' = '  -0.07%  '
    'E8' = '  -1.83%  '
    'E9' = '  -1.99%  '
    'E10' = '  -1.74%  '
    'B11' = 'Polygon'
    'C11' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E11' = '  -5.54%  '
    'B12' = 'BinanceUSD'
    'C12' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'E12' = '  +0.61%  '
    'E13' = '  -3.76%  '
    'E14' = '  -2.33%  '
    'E15' = '  -4.39%  '
    'E16' = '  -3.89%  '
    'D17' = '1.621.83'
    'E17' = '  -0.83%  '
    'E18' = '  -1.58%  '
    'E19' = '  -0.52%  '
    'E20' = '  -3.47%  '
    'E22' = '  -3.24%  '
    'D23' = '23.311.21'
    'E23' = '  -1.08%  '
    'E24' = '  -2.67%  '
    'E25' = '  +0.65%  '
    'E26' = '  +0.64%  '
    'E27' = '  -1.88%  '
    'E28' = '  -0.94%  '
    'E29' = '  -0.43%  '
    'E30' = '  -1.61%  '
    'E31' = '  -5.41%  '
    'D32' = '1.799.86'
    'E32' = '  -0.67%  '
    'E33' = '  -2.51%  '
    'E34' = '  +3.91%  '
    'E35' = '  -3.96%  '
    'E36' = '  -1.99%  '
    'E37' = '  -1.43%  '
    'E38' = '  -0.56%  '
    'E39' = '  -3.90%  '
    'E40' = '  -5.50%  '
    'E41' = '  -3.99%  '
    'E42' = '  -2.70%  '
    'E43' = '  -1.54%  '
    'E44' = '  -4.59%  '
    'E45' = '  +0.58%  '
    'E46' = '  -3.82%  '
    'E47' = '  -3.03%  '
    'E48' = '  -1.29%  '
    'E49' = '  -1.06%  '
    'E50' = '  -2.43%  '
    'E51' = '  -5.67%  '
}
foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}
